# Backup QR Scanner data update
# 1) Rename the worksheet from "Session" to "Pharmacology"
# 2) Append two new scan log rows (21 and 22)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Pharmacology"

# Ensure new rows keep their values as plain text, matching the existing
# "number stored as text" style used throughout the log (t="str").
$newRange = $ws.Range("A21:F22")
$newRange.NumberFormat = "@"

# Row 21
$ws.Cells.Item(21, 1).Value = "244939"
$ws.Cells.Item(21, 2).Value = "Pharmacology"
$ws.Cells.Item(21, 3).Value = "11/12/2025"
$ws.Cells.Item(21, 4).Value = "08:46:09"
$ws.Cells.Item(21, 5).Value = "Scan"
$ws.Cells.Item(21, 6).Value = "nancy.abdelshafy@med.asu.edu.eg"

# Row 22
$ws.Cells.Item(22, 1).Value = "244826"
$ws.Cells.Item(22, 2).Value = "Pharmacology"
$ws.Cells.Item(22, 3).Value = "11/12/2025"
$ws.Cells.Item(22, 4).Value = "08:51:07"
$ws.Cells.Item(22, 5).Value = "Scan"
$ws.Cells.Item(22, 6).Value = "nancy.abdelshafy@med.asu.edu.eg"

# Extend the "number stored as text" ignored-error range to cover the
# newly added rows, matching the updated dimension (A1:F22).
$fullRange = $ws.Range("A1:F22")
$errs = $fullRange.Errors
$errs.Item(6).Ignore = $true
